# Applies the crypto price/volume updates from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Writing a numeric-looking string via .Value lets Excel silently
    # coerce it to a Number (dropping things like trailing zeros).
    # Route it through a text-literal formula, then paste-special as
    # values (xlPasteValues = -4163) to freeze it back to a plain
    # string cell without leaving a numeric format/style behind.
    $Range.Formula = '="' + $Text + '"'
    $Range.Copy() | Out-Null
    $Range.PasteSpecial(-4163) | Out-Null
}

$ws.Range("D2").Value = "28.293.67"
$ws.Range("D3").Value = "1.817.93"
Set-TextValue $ws.Range("D4") "1.004"
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws.Range("D5") "328.13"
$ws.Range("E5").Value = "  +1.68%  "
Set-TextValue $ws.Range("D6") "1.002"
$ws.Range("E6").Value = "  +0.16%  "
Set-TextValue $ws.Range("D7") "0.4347"
$ws.Range("E7").Value = "  +3.26%  "
Set-TextValue $ws.Range("D8") "0.3671"
$ws.Range("E8").Value = "  +2.83%  "
Set-TextValue $ws.Range("D9") "44.92"
$ws.Range("E9").Value = "  -1.47%  "
Set-TextValue $ws.Range("D10") "0.07694"
$ws.Range("E10").Value = "  +3.86%  "
Set-TextValue $ws.Range("D11") "1.145"
$ws.Range("E11").Value = "  +3.01%  "
Set-TextValue $ws.Range("D12") "1.002"
$ws.Range("E12").Value = "  +0.07%  "
Set-TextValue $ws.Range("D13") "22.14"
$ws.Range("E13").Value = "  +3.32%  "
Set-TextValue $ws.Range("D14") "6.306"
$ws.Range("E14").Value = "  +3.40%  "
Set-TextValue $ws.Range("D15") "7.542"
$ws.Range("E15").Value = "  +5.19%  "
$ws.Range("D16").Value = "1.833.43"
$ws.Range("E16").Value = "  +5.33%  "
Set-TextValue $ws.Range("D17") "93.70"
$ws.Range("E17").Value = "  +7.14%  "
Set-TextValue $ws.Range("D18") "0.00001084"
$ws.Range("E18").Value = "  +1.94%  "
Set-TextValue $ws.Range("D19") "0.06566"
$ws.Range("E19").Value = "  +6.53%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E21").Value = "  +3.99%  "
Set-TextValue $ws.Range("D22") "6.269"
$ws.Range("E22").Value = "  +2.76%  "
$ws.Range("D23").Value = "28.329.29"
$ws.Range("E23").Value = "  +3.01%  "
Set-TextValue $ws.Range("D24") "11.62"
$ws.Range("E24").Value = "  +0.39%  "
Set-TextValue $ws.Range("D25") "2.067"
$ws.Range("E25").Value = "  -10.98%  "
Set-TextValue $ws.Range("D26") "162.78"
$ws.Range("E26").Value = "  +6.84%  "
Set-TextValue $ws.Range("D27") "20.71"
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("D28").Value = "2.040.89"
$ws.Range("E28").Value = "  +5.27%  "
Set-TextValue $ws.Range("D29") "2.300"
$ws.Range("E29").Value = "  -2.68%  "
Set-TextValue $ws.Range("D30") "129.01"
$ws.Range("E30").Value = "  +2.46%  "
Set-TextValue $ws.Range("D31") "1.230"
$ws.Range("E31").Value = "  +2.21%  "
Set-TextValue $ws.Range("D32") "5.977"
$ws.Range("E32").Value = "  +5.15%  "
Set-TextValue $ws.Range("D33") "0.09189"
$ws.Range("E33").Value = "  +0.53%  "
Set-TextValue $ws.Range("D34") "3.480"
$ws.Range("E34").Value = "  -5.62%  "
Set-TextValue $ws.Range("D35") "13.00"
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("E36").Value = "  +3.05%  "
Set-TextValue $ws.Range("D37") "0.2181"
$ws.Range("E37").Value = "  +2.63%  "
Set-TextValue $ws.Range("D38") "5.207"
$ws.Range("E38").Value = "  +2.58%  "
Set-TextValue $ws.Range("D39") "0.6584"
$ws.Range("E39").Value = "  +3.16%  "
Set-TextValue $ws.Range("D40") "0.06208"
$ws.Range("E40").Value = "  +2.26%  "
Set-TextValue $ws.Range("D41") "1.196"
Set-TextValue $ws.Range("D42") "8.142"
$ws.Range("E42").Value = "  +3.20%  "
Set-TextValue $ws.Range("D43") "1.438"
$ws.Range("E43").Value = "  +1.10%  "
Set-TextValue $ws.Range("D44") "1.001"
$ws.Range("E44").Value = "  +0.11%  "
Set-TextValue $ws.Range("D45") "13.95"
$ws.Range("E45").Value = "  +1.44%  "
Set-TextValue $ws.Range("D46") "0.6121"
$ws.Range("E46").Value = "  +4.63%  "
$ws.Range("E47").Value = "  +1.09%  "
Set-TextValue $ws.Range("D48") "2.028"
$ws.Range("E48").Value = "  +4.09%  "
Set-TextValue $ws.Range("D49") "126.00"
$ws.Range("E49").Value = "  +0.79%  "
Set-TextValue $ws.Range("D50") "1.159"
$ws.Range("E50").Value = "  +3.38%  "
Set-TextValue $ws.Range("D51") "0.07011"
$ws.Range("E51").Value = "  +2.41%  "

$excel.CutCopyMode = 0
